$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = "Google Cloud Platform DevOps Engneer (Telecom domain & Visa independent only)"
$ws.Range("B31").Value = "https://www.dice.com/job-detail/88da6340-a700-441c-be19-26365a5582fe"
$ws.Range("C31").Value = "Remote or San Jose, California"
$ws.Range("D31").Value = "Contract, Third Party"
$ws.Range("E31").Value = "Depends on Experience"
$ws.Range("F31").Value = "Zeforge LLC"
